# Update the "想去人数" (interest count) values in column F on the
# "展览" (Exhibitions) and "全部类型" (All Types) sheets to reflect the
# latest scrape of generated output (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 147
$ws1.Range("F3").Value  = 1337
$ws1.Range("F4").Value  = 1144
$ws1.Range("F5").Value  = 1033
$ws1.Range("F6").Value  = 1813
$ws1.Range("F7").Value  = 576
$ws1.Range("F8").Value  = 1213
$ws1.Range("F12").Value = 305
$ws1.Range("F13").Value = 80
$ws1.Range("F15").Value = 704
$ws1.Range("F16").Value = 184
$ws1.Range("F21").Value = 163
$ws1.Range("F22").Value = 679
$ws1.Range("F23").Value = 46
$ws1.Range("F24").Value = 653
$ws1.Range("F25").Value = 164
$ws1.Range("F27").Value = 883
$ws1.Range("F28").Value = 323
$ws1.Range("F29").Value = 166
$ws1.Range("F30").Value = 46
$ws1.Range("F31").Value = 279

# Sheet "全部类型" (sheet4): row -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 147
$ws4.Range("F4").Value  = 1337
$ws4.Range("F5").Value  = 1144
$ws4.Range("F6").Value  = 1033
$ws4.Range("F7").Value  = 1813
$ws4.Range("F8").Value  = 576
$ws4.Range("F9").Value  = 1213
$ws4.Range("F14").Value = 305
$ws4.Range("F15").Value = 80
$ws4.Range("F17").Value = 704
$ws4.Range("F18").Value = 184
$ws4.Range("F29").Value = 163
$ws4.Range("F30").Value = 679
$ws4.Range("F31").Value = 46
$ws4.Range("F32").Value = 653
$ws4.Range("F33").Value = 164
$ws4.Range("F35").Value = 883
$ws4.Range("F36").Value = 323
$ws4.Range("F39").Value = 166
$ws4.Range("F40").Value = 46
$ws4.Range("F41").Value = 279
